$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.455.96'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '1.624.77'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("E6").Value = '  +1.40%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("E9").Value = '  -1.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0841'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").Value = '1.851.84'
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D13").Value = '1.634.84'
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.34%  '
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.88'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.01%  '
$ws.Range("D17").Value = '26.502.34'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '213.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.47%  '
$ws.Range("E21").Value = '  -0.68%  '
$ws.Range("E22").Value = '  +1.54%  '
$ws.Range("E23").Value = '  -1.26%  '
$ws.Range("E24").Value = '  +7.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("E27").Value = '  -0.67%  '
$ws.Range("E28").Value = '  +1.05%  '
$ws.Range("E29").Value = '  +0.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0506'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.56%  '
$ws.Range("E31").Value = '  -1.22%  '
$ws.Range("E32").Value = '  +2.70%  '
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("D35").Value = '1.227.51'
$ws.Range("E35").Value = '  +5.19%  '
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0172'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.35%  '
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.791'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.503'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.791'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.33'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.85%  '
$ws.Range("D44").Value = '1.760.93'
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.66'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0509'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("E49").Value = '  -0.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.46'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.85%  '
$ws.Range("E51").Value = '  +0.17%  '
